# add more tag and keywords
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A36").Value = "มีคำถาม"
$ws.Range("B36").Value = "ถามได้เลยค้าบบ"

$ws.Range("A37").Value = "มีคำถาม"
$ws.Range("B37").Value = "ถามมาได้เลยค้าบบ"

$ws.Range("I49").Select()
